$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 16 de Agosto de 2020 a las 02:02"

# Re-rank Argentina ahead of Pakistan (rows 17-18)
$ws.Range("A17").Value = "Argentina"
$ws.Range("A18").Value = "Pakistan"

# Re-rank Surinam ahead of Sri Lanka/Eslovaquia (rows 122-124)
$ws.Range("A122").Value = "Surinam"
$ws.Range("A123").Value = "Sri Lanka"
$ws.Range("A124").Value = "Eslovaquia"

# Update numeric statistics for all affected rows
# Row 4
$ws.Range("B4").Value = 5529289
$ws.Range("C4").Value = 53023
$ws.Range("D4").Value = 2898566
$ws.Range("E4").Value = 2458138
$ws.Range("G4").Value = 1050
$ws.Range("H4").Value = 172585

# Row 5
$ws.Range("B5").Value = 3317832
$ws.Range("C5").Value = 38937
$ws.Range("E5").Value = 806263
$ws.Range("G5").Value = 726
$ws.Range("H5").Value = 107297

# Row 17
$ws.Range("B17").Value = 289100
$ws.Range("C17").Value = 6663
$ws.Range("D17").Value = 205697
$ws.Range("E17").Value = 77766
$ws.Range("G17").Value = 110
$ws.Range("H17").Value = 5637

# Row 18
$ws.Range("B18").Value = 288047
$ws.Range("C18").Value = 747
$ws.Range("D18").Value = 265624
$ws.Range("E18").Value = 16261
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 6162

# Row 74
$ws.Range("B74").Value = 19891
$ws.Range("C74").Value = 198
$ws.Range("D74").Value = 13763
$ws.Range("E74").Value = 5733

# Row 81
$ws.Range("B81").Value = 14333
$ws.Range("C81").Value = 90
$ws.Range("D81").Value = 9161
$ws.Range("E81").Value = 4677
$ws.Range("G81").Value = 3
$ws.Range("H81").Value = 495

# Row 98
$ws.Range("B98").Value = 7439
$ws.Range("C98").Value = 34
$ws.Range("E98").Value = 816
$ws.Range("G98").Value = 1
$ws.Range("H98").Value = 123

# Row 113
$ws.Range("B113").Value = 3960
$ws.Range("C113").Value = 30
$ws.Range("D113").Value = 2830
$ws.Range("E113").Value = 1055
$ws.Range("G113").Value = 2
$ws.Range("H113").Value = 75

# Row 122
$ws.Range("B122").Value = 2961
$ws.Range("C122").Value = 123
$ws.Range("D122").Value = 1990
$ws.Range("E122").Value = 929
$ws.Range("G122").Value = 1
$ws.Range("H122").Value = 42

# Row 123
$ws.Range("B123").Value = 2890
$ws.Range("C123").Value = 4
$ws.Range("D123").Value = 2666
$ws.Range("E123").Value = 213
$ws.Range("H123").Value = 11

# Row 124
$ws.Range("B124").Value = 2855
$ws.Range("C124").Value = 54
$ws.Range("D124").Value = 1969
$ws.Range("E124").Value = 855
$ws.Range("H124").Value = 31

# Row 132
$ws.Range("B132").Value = 2117
$ws.Range("C132").Value = 29
$ws.Range("E132").Value = 1069
$ws.Range("G132").Value = 4
$ws.Range("H132").Value = 33

# Row 143
$ws.Range("B143").Value = 1434
$ws.Range("C143").Value = 13
$ws.Range("D143").Value = 1194
$ws.Range("E143").Value = 202

# Row 159
$ws.Range("B159").Value = 951
$ws.Range("C159").Value = 21
$ws.Range("E159").Value = 481

# Row 161
$ws.Range("D161").Value = 815
$ws.Range("E161").Value = 55

